$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the hours-per-person grid (columns B..G = Re, Ad, An, Pr, Pg, Ve).
# Blank cells in the original sheet become literal "-" text cells.
# Each entry is: row, Re, Ad, An, Pr, Pg, Ve
$data = @(
    ,(2, 2,   "-", 1,   "-", "-", 2)
    ,(3, "-", "-", 3,   "-", "-", 2)
    ,(4, "-", 1,   4,   "-", "-", "-")
    ,(5, "-", "-", 3,   "-", "-", 2)
    ,(6, 2,   "-", 3,   "-", "-", "-")
    ,(7, "-", 2,   "-", "-", "-", 3)
    ,(8, "-", "-", 4,   "-", "-", 1)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]  # B - Re
    $ws.Cells.Item($row, 3).Value = $entry[2]  # C - Ad
    $ws.Cells.Item($row, 4).Value = $entry[3]  # D - An
    $ws.Cells.Item($row, 5).Value = $entry[4]  # E - Pr
    $ws.Cells.Item($row, 6).Value = $entry[5]  # F - Pg
    $ws.Cells.Item($row, 7).Value = $entry[6]  # G - Ve
}

# Reposition / resize the chart (Grafico 1) to its new anchor.
$co = $ws.ChartObjects(1)
$co.Left = 706.8397650098425
$co.Top = 19.68748031496063
$co.Width = 748.4375
$co.Height = 287.28748031496065

# Leave the active selection on A9, matching the saved view state.
$null = $ws.Range("A9").Select()
